$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a brand new row below row 7 (pushing the old totals/footer rows
#    down) and clone row 7's formatting onto it, so both sale rows share the
#    same look (borders, fills, number formats, merges).
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).Insert()
$ws.Range("A7:Q7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5

$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# ---------------------------------------------------------------------------
# 2) Fill in row 7 with the first sold item (was a blank placeholder row).
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 1

$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("C7").Value = "EREC 100MG 12 F.C. TABLETS"

$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "6:8"

$ws.Range("L7").Value = "'1"

$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("N7").Value = "144.00"

$ws.Range("P7").Value = "'155.5200"

$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1:1"

# ---------------------------------------------------------------------------
# 3) Fill in the new row 8 with the second sold item.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 2

$ws.Range("C8:G8").NumberFormat = "@"
$ws.Range("C8").Value = "NETLOOK 10MG 20 SOFT GELATIN CAPS."

$ws.Range("H8:K8").NumberFormat = "@"
$ws.Range("H8").Value = "1:0"

$ws.Range("L8").Value = "'1"

$ws.Range("N8:O8").NumberFormat = "@"
$ws.Range("N8").Value = "150.00"

$ws.Range("P8").Value = "'150.0000"

$ws.Range("Q8").NumberFormat = "@"
$ws.Range("Q8").Value = "1:0"

# ---------------------------------------------------------------------------
# 4) The (now shifted) totals row gets the sum of the sale prices.
# ---------------------------------------------------------------------------
$ws.Range("P9").Value = 305.51999999999998

# ---------------------------------------------------------------------------
# 5) The footer row (now row 10) keeps its text but the timestamp advances
#    one minute, matching the re-export.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Wednesday, 30 July, 2025 12:32 AM"
